$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the dummy sample rows (rows 2-5); only the header row survives.
$ws.Rows("2:5").Delete()

# Re-label the header row with the new column set/order. A1:F1 already
# carry the bordered/bold header style, so overwriting their values keeps
# that formatting; new columns G1:I1 get the same style copied onto them.
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "member_id"
$ws.Range("D1").Value = "member_name"
$ws.Range("E1").Value = "inventory_id"
$ws.Range("F1").Value = "item_name"
$ws.Range("G1").Value = "quantity"
$ws.Range("H1").Value = "total_amount"
$ws.Range("I1").Value = "payment_method"

# Carry the header styling (bold, border, centered) onto the new columns.
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
